# Adds a fourth "Instant APs" Q&A block to the mock_data table (rows 22-26)
# and grows the table / sheet dimension / view selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block (mirrors the existing "#, Question, Data Source, Relevance"
#     pattern used by the other blocks, with a blank separator row at 21) ---
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Instant APs"
$ws.Range("C22").Value = "iaps.htm"
$ws.Range("D22").Value = 2

$ws.Range("B23").Value = "What are instant APs?"
$ws.Range("B24").Value = "define Instant APs"
$ws.Range("B25").Value = "What are IAPs?"
$ws.Range("B26").Value = "Define IAPs"

# --- Grow the table / autofilter range so the new rows are included ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D26"))

# --- Update the view: scroll down a bit and select B28 (matches the saved
#     sheetView's topLeftCell/selection in the edited workbook) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B28").Select()
